$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text in B1: remove the leading space from "ErrorDescription"
$ws.Range("B1").Value = "ErrorDescription"

# Update the active selection to G8 (matches the sheetView selection in the target)
$ws.Range("G8").Select()
